$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 30 new rows above the existing 6 rows of data, shifting them down to rows 31-36
$ws.Range("A1:A30").EntireRow.Insert()

$names = @(
    'Gigabyte GeForce RTX 3060 EAGLE 12GB OC 2.0 LHR',
    'Gigabyte GeForce RTX 3060 Gaming 12GB OC 2.0 LHR',
    'ASUS GeForce RTX 3060 DUAL 12GB OC V2',
    'ASUS GeForce RTX 3060 TUF 12GB OC V2 LHR',
    'ASUS GeForce RTX 3060 DUAL 8GB OC WHITE',
    'Gigabyte GeForce RTX 3060 WindForce 12GB OC',
    'ASUS GeForce RTX 3060 DUAL 8GB OC',
    'Gigabyte GeForce RTX 3060 VISION 12GB OC 2.0 LHR',
    'MSI GeForce RTX 3060 VENTUS 2X 12GB OC LHR',
    'Gainward GeForce RTX 3060 Ghost 12GB',
    'ASUS GeForce RTX 3060 PHOENIX 12GB V2 LHR',
    'ZOTAC Gaming GeForce RTX 3060 Twin Edge OC 12GB',
    'Inno3D GeForce RTX 3060 Twin X2 LHR 12GB GDDR6',
    'Gigabyte GeForce RTX 3060 Gaming 8GB OC',
    'Palit GeForce RTX 3060 Dual 12GB',
    'ZOTAC Gaming GeForce RTX 3060 Twin Edge 12GB',
    'Palit GeForce RTX 3060 Dual 12GB OC',
    'Inno3D GeForce RTX 3060 Twin X2 OC 8GB GDDR6',
    'MSI GeForce RTX 3060 GAMING X 12GB LHR',
    'MSI GeForce RTX 3060 VENTUS 3X 12GB OC LHR',
    'MSI GeForce RTX 3060 VENTUS 2X 8GB OC LHR',
    'MSI GeForce RTX 3060 AERO ITX 12GB OC LHR',
    'ZOTAC Gaming GeForce RTX 3060 AMP White Edition 12GB',
    'Inno3D GeForce RTX 3060 Twin X2 OC LHR 12GB GDDR6',
    'Gainward GeForce RTX 3060 Ghost 12GB OC',
    'Gigabyte GeForce RTX 3060 AORUS ELITE 12GB 2.0 LHR',
    'Gainward GeForce RTX 3060 Pegasus OC 12GB',
    'Gainward GeForce RTX 3060 Pegasus 12GB',
    'ZOTAC Gaming GeForce RTX 3060 Twin Edge 8GB',
    'MSI GeForce RTX 3060 GAMING Z 12GB LHR',
    'PNY GeForce RTX 3060 12GB XLR8 Gaming REVEL EPIC-X RGB Single Fan',
    'Palit GeForce RTX 3060 StormX',
    'PNY GeForce RTX 3060 12GB Uprising Dual Fan',
    'PNY GeForce RTX 3060 VERTO 8GB Dual Fan Edition',
    'Inno3D GeForce RTX 3060 Twin X2 8GB GDDR6',
    'PNY GeForce RTX 3060 12GB XLR8 Gaming REVEL EPIC-X RGB Dual Fan',
)

$prices = @(
    1519,
    1499,
    1629,
    1829,
    1429,
    1569,
    1409,
    1839,
    1589,
    1379,
    1349,
    1469,
    1399,
    1609,
    1399,
    1459,
    1449,
    1409,
    1659,
    1599,
    1449,
    2049,
    1579,
    1589,
    1779,
    2399,
    3699,
    2199,
    1439,
    1809,
    1779,
    3699,
    4749,
    1459,
    1679,
    2049,
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $prices[$i]
}

